$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 492 (shifts existing rows 492:617 down to 493:618)
$ws.Rows.Item(492).Insert()

# Populate the newly inserted row 492 with the new data point
$ws.Range("A492").Value2 = 3
$ws.Range("B492").Value2 = "Femacal de La Calera"
$ws.Range("C492").Value2 = "Coquimbo"
$ws.Range("D492").Value2 = 45204
$ws.Range("E492").Value2 = 5
$ws.Range("F492").Value2 = 100112012
$ws.Range("G492").Value2 = "Espinaca"
$ws.Range("H492").Value2 = "Sin especificar"
$ws.Range("I492").Value2 = "Primera"
$ws.Range("J492").Value2 = 90
$ws.Range("K492").Value2 = 4000
$ws.Range("L492").Value2 = 4000
$ws.Range("M492").Value2 = 4000
$ws.Range("N492").Value2 = "`$/docena de atados (3 kilos)"
$ws.Range("O492").Value2 = "Provincia de Quillota"
$ws.Range("P492").Value2 = 1333
$ws.Range("Q492").Value2 = 3
$ws.Range("R492").Value2 = "Hortaliza"
